$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (logistic_embeddings)
$ws.Range("C5").Value = 0.222
$ws.Range("D5").Value = 0.38
$ws.Range("E5").Value = 0.432
$ws.Range("F5").Value = 0.46
$ws.Range("G5").Value = 0.497
$ws.Range("H5").Value = 0.515

# Row 7 (classical-best-embeddings -> classical-best-embed)
$ws.Range("A7").Value = "classical-best-embed"
$ws.Range("C7").Value = 0.222
$ws.Range("D7").Value = 0.38
$ws.Range("E7").Value = 0.432

# Row 8 (BERT-base)
$ws.Range("C8").Value = 0.202
$ws.Range("D8").Value = 0.477
$ws.Range("E8").Value = 0.514
$ws.Range("F8").Value = 0.576
$ws.Range("G8").Value = 0.618
$ws.Range("H8").Value = 0.619

# Row 9 (BERT-base-nli)
$ws.Range("B9").Value = 0.358
$ws.Range("C9").Value = 0.49
$ws.Range("D9").Value = 0.574
$ws.Range("E9").Value = 0.588
$ws.Range("F9").Value = 0.635
$ws.Range("G9").Value = 0.656
$ws.Range("H9").Value = 0.658
